$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 225.5
$ws.Range("I4").Value = 225.5
$ws.Range("K4").Value = 225.5
$ws.Range("M4").Value = -111.5

$ws.Range("H107").Value = 2273.1702
$ws.Range("I107").Value = 1941.5938
$ws.Range("J107").Value = 2980.5334
$ws.Range("K107").Value = 1941.5938
$ws.Range("L107").Value = 2980.5334
$ws.Range("M107").Value = -21.5938000000001
$ws.Range("N107").Value = -6820.5334

$ws.Range("H137").Value = 3336.4
$ws.Range("I137").Value = 4837.4814
$ws.Range("K137").Value = 14512.4442
$ws.Range("M137").Value = -11962.4442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 83334.2
$ws.Range("J24").Value = 83334.2
$ws.Range("L24").Value = 83334.2
$ws.Range("N24").Value = -84082.2

$ws.Range("H32").Value = 3146.4023
$ws.Range("I32").Value = 1974.3422
$ws.Range("J32").Value = 17992.5
$ws.Range("K32").Value = 1974.3422
$ws.Range("L32").Value = 17992.5
$ws.Range("M32").Value = -1687.3422
$ws.Range("N32").Value = -18566.5

$ws.Range("H74").Value = 22181.928
$ws.Range("I74").Value = 32949.5
$ws.Range("J74").Value = 17874.9
$ws.Range("K74").Value = 32949.5
$ws.Range("L74").Value = 17874.9
$ws.Range("M74").Value = -32075.5
$ws.Range("N74").Value = -19622.9

$ws.Range("H77").Value = 22181.928
$ws.Range("I77").Value = 32949.5
$ws.Range("J77").Value = 17874.9
$ws.Range("K77").Value = 164747.5
$ws.Range("L77").Value = 89374.5
$ws.Range("M77").Value = -160379.5
$ws.Range("N77").Value = -98110.5

$ws.Range("H100").Value = 83334.2
$ws.Range("J100").Value = 83334.2
$ws.Range("L100").Value = 83334.2
$ws.Range("N100").Value = -85498.2

$ws.Range("H102").Value = 23836450
$ws.Range("J102").Value = 127147.5
$ws.Range("L102").Value = 127147.5
$ws.Range("N102").Value = -130391.5

$ws.Range("H110").Value = 2097.75
$ws.Range("I110").Value = 1754.5714
$ws.Range("J110").Value = 4500
$ws.Range("K110").Value = 1754.5714
$ws.Range("L110").Value = 4500
$ws.Range("M110").Value = 290.4286
$ws.Range("N110").Value = -8590

$ws.Range("H125").Value = 38901.08
$ws.Range("J125").Value = 38901.08
$ws.Range("L125").Value = 38901.08
$ws.Range("N125").Value = -48741.08

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4727.5713
$ws.Range("I86").Value = 5231.3335
$ws.Range("J86").Value = 4349.75
$ws.Range("K86").Value = 5231.3335
$ws.Range("L86").Value = 4349.75
$ws.Range("M86").Value = -4108.3335
$ws.Range("N86").Value = -6595.75

$ws.Range("H89").Value = 4727.5713
$ws.Range("I89").Value = 5231.3335
$ws.Range("J89").Value = 4349.75
$ws.Range("K89").Value = 26156.6675
$ws.Range("L89").Value = 21748.75
$ws.Range("M89").Value = -20540.6675
$ws.Range("N89").Value = -32980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3085.5833
$ws.Range("I58").Value = 3279.375
$ws.Range("K58").Value = 3279.375
$ws.Range("M58").Value = -3076.375

$ws.Range("H132").Value = 1862.2258
$ws.Range("I132").Value = 1645.8276
$ws.Range("K132").Value = 4937.4828
$ws.Range("M132").Value = -2407.4828

$ws.Range("H136").Value = 3085.5833
$ws.Range("I136").Value = 3279.375
$ws.Range("K136").Value = 9838.125
$ws.Range("M136").Value = -7288.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 174.41667
$ws.Range("I40").Value = 157.3
$ws.Range("J40").Value = 260
$ws.Range("K40").Value = 629.2
$ws.Range("L40").Value = 1040
$ws.Range("M40").Value = -560.2
$ws.Range("N40").Value = -1178

$ws.Range("H120").Value = 13816.833
$ws.Range("I120").Value = 11633.667
$ws.Range("K120").Value = 34901.001
$ws.Range("M120").Value = -30063.001

$ws.Range("H129").Value = 2764.6943
$ws.Range("I129").Value = 3015.3333
$ws.Range("J129").Value = 2639.375
$ws.Range("K129").Value = 9045.999899999999
$ws.Range("L129").Value = 7918.125
$ws.Range("M129").Value = -4045.999899999999
$ws.Range("N129").Value = -17918.125

$ws.Range("H131").Value = 3448
$ws.Range("I131").Value = 3674.5
$ws.Range("K131").Value = 11023.5
$ws.Range("M131").Value = -5983.5

$ws.Range("H134").Value = 4020.7144
$ws.Range("I134").Value = 1430
$ws.Range("K134").Value = 4290
$ws.Range("M134").Value = 780

$ws.Range("H140").Value = 1449.3846
$ws.Range("I140").Value = 1195.1666
$ws.Range("K140").Value = 3585.4998
$ws.Range("M140").Value = 1594.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 500063
$ws.Range("I70").Value = 500063
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 500063
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -499793
$ws.Range("N70").Value = ""

$ws.Range("H73").Value = 500063
$ws.Range("I73").Value = 500063
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 500063
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -499127
$ws.Range("N73").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58831124
$ws.Range("J7").Value = 8926.25
$ws.Range("L7").Value = 8926.25
$ws.Range("N7").Value = -9150.25

$ws.Range("H22").Value = 21947.271
$ws.Range("I22").Value = 251083
$ws.Range("J22").Value = 1116.75
$ws.Range("K22").Value = 251083
$ws.Range("L22").Value = 1116.75
$ws.Range("M22").Value = -250788
$ws.Range("N22").Value = -1706.75

$ws.Range("H27").Value = 21947.271
$ws.Range("I27").Value = 251083
$ws.Range("J27").Value = 1116.75
$ws.Range("K27").Value = 251083
$ws.Range("L27").Value = 1116.75
$ws.Range("M27").Value = -250976
$ws.Range("N27").Value = -1330.75

$ws.Range("H55").Value = 1776
$ws.Range("I55").Value = 611.7
$ws.Range("J55").Value = 3231.375
$ws.Range("K55").Value = 611.7
$ws.Range("L55").Value = 3231.375
$ws.Range("M55").Value = -438.7
$ws.Range("N55").Value = -3577.375

$ws.Range("H61").Value = 1899.5
$ws.Range("I61").Value = 1899.5
$ws.Range("K61").Value = 1899.5
$ws.Range("M61").Value = -1697.5

$ws.Range("H82").Value = 90914730
$ws.Range("I82").Value = 200001360
$ws.Range("J82").Value = 9202.833000000001
$ws.Range("K82").Value = 200001360
$ws.Range("L82").Value = 9202.833000000001
$ws.Range("M82").Value = -200000999
$ws.Range("N82").Value = -9924.833000000001

$ws.Range("H85").Value = 90914730
$ws.Range("I85").Value = 200001360
$ws.Range("J85").Value = 9202.833000000001
$ws.Range("K85").Value = 200001360
$ws.Range("L85").Value = 9202.833000000001
$ws.Range("M85").Value = -200000112
$ws.Range("N85").Value = -11698.833

$ws.Range("H100").Value = 5299.5835
$ws.Range("I100").Value = 3942.4285
$ws.Range("J100").Value = 7199.6
$ws.Range("K100").Value = 3942.4285
$ws.Range("L100").Value = 7199.6
$ws.Range("M100").Value = -3401.4285
$ws.Range("N100").Value = -8281.6

$ws.Range("H113").Value = 1899.5
$ws.Range("I113").Value = 1899.5
$ws.Range("K113").Value = 1899.5
$ws.Range("M113").Value = 270.5

$ws.Range("H126").Value = 58831124
$ws.Range("J126").Value = 8926.25
$ws.Range("L126").Value = 26778.75
$ws.Range("N126").Value = -31718.75

$ws.Range("H136").Value = 5617.0713
$ws.Range("I136").Value = 5967.263
$ws.Range("K136").Value = 17901.789
$ws.Range("M136").Value = -15351.789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

$ws.Range("H132").Value = 1635.1765
$ws.Range("I132").Value = 1473.2667
$ws.Range("K132").Value = 4419.800099999999
$ws.Range("M132").Value = -1889.800099999999
